$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_LB")

# Rename sheet tab from Collection_LB to CRF_LB.
# This automatically updates the sheet name references used by defined
# names (e.g. the hidden _FilterDatabase name) and document properties.
$ws.Name = "CRF_LB"

# Normalize the style of the previously-duplicate style (index 13) cells
# to the equivalent style (index 11): left-aligned, wrap text, no border.
$cells = @("R892", "R904", "R916", "R928", "R940", "R950")
foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = -4131   # xlLeft
    $rng.WrapText = $true
    $rng.Borders.LineStyle = -4142     # xlLineStyleNone
}
